$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- 1. Update existing "prices" sheet -------------------------------------------------
# C2 (Paris lowest price) 54 -> 76
$ws1.Range("C2").Value = 76

# New row 11: Bali / DPS / 501
$ws1.Range("A11").Value = "Bali"
$ws1.Range("B11").Value = "DPS"
$ws1.Range("C11").Value = 501

$row11 = $ws1.Range("A11:C11")
$row11.Font.ThemeColor = 1

# --- 2. Add the new "users" sheet ------------------------------------------------------
# Copy "prices" so the new sheet inherits the same sheet-level defaults
# (outline props, default row height / column width, etc.), then wipe its content.
$ws1.Copy($null, $ws1)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "users"
$ws2.Cells.Clear()

# Header row
$ws2.Range("A1").Value = "First Name"
$ws2.Range("B1").Value = "Last Name"
$ws2.Range("C1").Value = "Email"

$headerRow = $ws2.Range("A1:Z1")
$headerRow.Font.ThemeColor = 1
$headerRow.Font.Bold = $true

# Data rows
$ws2.Range("A2").Value = "Name1"
$ws2.Range("B2").Value = "LastName1"
$ws2.Range("C2").Value = "email1@gmail.com"

$ws2.Range("A3").Value = "Name2"
$ws2.Range("B3").Value = "LastName2"
$ws2.Range("C3").Value = "email2@gmail.com"

$dataRows = $ws2.Range("A2:C3")
$dataRows.Font.ThemeColor = 1

# Column C is wider to fit the email addresses
$ws2.Columns.Item(3).ColumnWidth = 19.3

Write-Output "done"
